# Add a new row (row 3) with the author's name, email and repo link,
# matching the existing "name / email / Repo Link" table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new row values ---------------------------------------------------
$ws.Range("A3").Value = "جبريل إسلام حنفي محمود"
$ws.Range("B3").Value = "gipreel1424@gmail.com"
$ws.Range("C3").Value = "https://github.com/gipreelislam/Metrics.git"

# --- hyperlinks for the email + repo link cells ------------------------
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:gipreel1424@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/gipreelislam/Metrics.git")

# --- match the formatting used by the row above (row 2: Hyperlink style)
# (Hyperlinks.Add above already nudges the style; re-apply the same
# look-and-feel as B2:C2 so the new cells are visually consistent.)
$ws.Range("B2:C2").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- leave the selection where the author ended up after data entry ----
[void]$ws.Range("B6").Select()

Write-Output "Added row 3 (name/email/repo link) to Sheet1"
